$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-02-08 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-02-09 Friday", 2) | Out-Null

$t = $d.Tables.Item(1)

# Update each data cell by position (row, col) to avoid collisions between old/new values
$t.Cell(1, 1).Range.Text = "94÷5=18, 4"  # was "50÷6=8, 2"
$t.Cell(1, 2).Range.Text = "70÷5=14, 0"  # was "54÷4=13, 2"
$t.Cell(1, 3).Range.Text = "11÷4=2, 3"  # was "92÷6=15, 2"
$t.Cell(1, 4).Range.Text = "51÷4=12, 3"  # was "32÷8=4, 0"
$t.Cell(1, 5).Range.Text = "51÷9=5, 6"  # was "50÷2=25, 0"
$t.Cell(5, 1).Range.Text = "93÷3=31, 0"  # was "76÷2=38, 0"
$t.Cell(5, 2).Range.Text = "37÷4=9, 1"  # was "11÷6=1, 5"
$t.Cell(5, 3).Range.Text = "23÷5=4, 3"  # was "75÷4=18, 3"
$t.Cell(5, 4).Range.Text = "84÷8=10, 4"  # was "28÷9=3, 1"
$t.Cell(5, 5).Range.Text = "49÷3=16, 1"  # was "15÷4=3, 3"
$t.Cell(9, 1).Range.Text = "76÷8=9, 4"  # was "41÷5=8, 1"
$t.Cell(9, 2).Range.Text = "76÷2=38, 0"  # was "74÷3=24, 2"
$t.Cell(9, 3).Range.Text = "23÷9=2, 5"  # was "92÷8=11, 4"
$t.Cell(9, 4).Range.Text = "79÷3=26, 1"  # was "56÷2=28, 0"
$t.Cell(9, 5).Range.Text = "31÷6=5, 1"  # was "79÷2=39, 1"
$t.Cell(13, 1).Range.Text = "12÷4=3, 0"  # was "19÷9=2, 1"
$t.Cell(13, 2).Range.Text = "25÷9=2, 7"  # was "62÷8=7, 6"
$t.Cell(13, 3).Range.Text = "99÷9=11, 0"  # was "82÷9=9, 1"
$t.Cell(13, 4).Range.Text = "26÷5=5, 1"  # was "33÷6=5, 3"
$t.Cell(13, 5).Range.Text = "73÷3=24, 1"  # was "84÷4=21, 0"
$t.Cell(17, 1).Range.Text = "11÷2=5, 1"  # was "98÷2=49, 0"
$t.Cell(17, 2).Range.Text = "44÷3=14, 2"  # was "27÷4=6, 3"
$t.Cell(17, 3).Range.Text = "41÷9=4, 5"  # was "51÷9=5, 6"
$t.Cell(17, 4).Range.Text = "64÷6=10, 4"  # was "50÷4=12, 2"
$t.Cell(17, 5).Range.Text = "43÷2=21, 1"  # was "71÷3=23, 2"
